$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(16, 3).Value2 = "45582159"
$ws.Cells.Item(16, 4).Value2 = "MAIDA CECILIA PIÑERES MARQUEZ"
$ws.Cells.Item(16, 5).Value2 = "1707"
$ws.Cells.Item(16, 6).Value2 = 29509
$ws.Cells.Item(16, 7).Value2 = 737717

$ws.Cells.Item(17, 3).Value2 = "45645018"
$ws.Cells.Item(17, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(17, 5).Value2 = "1707"
$ws.Cells.Item(17, 6).Value2 = 29509
$ws.Cells.Item(17, 7).Value2 = 781242

$ws.Cells.Item(18, 3).Value2 = "45582159"
$ws.Cells.Item(18, 4).Value2 = "MAIDA CECILIA PIÑERES MARQUEZ"
$ws.Cells.Item(18, 5).Value2 = "1708"
$ws.Cells.Item(18, 6).Value2 = 29509
$ws.Cells.Item(18, 7).Value2 = 737717

$ws.Cells.Item(19, 3).Value2 = "45645018"
$ws.Cells.Item(19, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(19, 5).Value2 = "1708"
$ws.Cells.Item(19, 6).Value2 = 29509
$ws.Cells.Item(19, 7).Value2 = 781242

$ws.Cells.Item(20, 3).Value2 = "45582159"
$ws.Cells.Item(20, 4).Value2 = "MAIDA CECILIA PIÑERES MARQUEZ"
$ws.Cells.Item(20, 5).Value2 = "1709"
$ws.Cells.Item(20, 6).Value2 = 29509
$ws.Cells.Item(20, 7).Value2 = 737717

$ws.Cells.Item(21, 3).Value2 = "45645018"
$ws.Cells.Item(21, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(21, 5).Value2 = "1709"
$ws.Cells.Item(21, 6).Value2 = 29509
$ws.Cells.Item(21, 7).Value2 = 781242

$ws.Cells.Item(22, 3).Value2 = "45582159"
$ws.Cells.Item(22, 4).Value2 = "MAIDA CECILIA PIÑERES MARQUEZ"
$ws.Cells.Item(22, 5).Value2 = "1710"
$ws.Cells.Item(22, 6).Value2 = 29509
$ws.Cells.Item(22, 7).Value2 = 737717

$ws.Cells.Item(23, 3).Value2 = "45645018"
$ws.Cells.Item(23, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(23, 5).Value2 = "1710"
$ws.Cells.Item(23, 6).Value2 = 29509
$ws.Cells.Item(23, 7).Value2 = 781242

$ws.Cells.Item(24, 3).Value2 = "45645018"
$ws.Cells.Item(24, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(24, 5).Value2 = "1711"
$ws.Cells.Item(24, 6).Value2 = 29509
$ws.Cells.Item(24, 7).Value2 = 781242

$ws.Cells.Item(25, 3).Value2 = "45645018"
$ws.Cells.Item(25, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(25, 5).Value2 = "1712"
$ws.Cells.Item(25, 6).Value2 = 29509
$ws.Cells.Item(25, 7).Value2 = 781242

$ws.Cells.Item(26, 3).Value2 = "45645018"
$ws.Cells.Item(26, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(26, 5).Value2 = "1801"
$ws.Cells.Item(26, 6).Value2 = 29509
$ws.Cells.Item(26, 7).Value2 = 781242

$ws.Cells.Item(27, 3).Value2 = "45645018"
$ws.Cells.Item(27, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(27, 5).Value2 = "1802"
$ws.Cells.Item(27, 6).Value2 = 29509
$ws.Cells.Item(27, 7).Value2 = 781242

$ws.Cells.Item(28, 3).Value2 = "45645018"
$ws.Cells.Item(28, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(28, 5).Value2 = "1803"
$ws.Cells.Item(28, 6).Value2 = 29509
$ws.Cells.Item(28, 7).Value2 = 781242

$ws.Cells.Item(29, 3).Value2 = "45645018"
$ws.Cells.Item(29, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(29, 5).Value2 = "1804"
$ws.Cells.Item(29, 6).Value2 = 29509
$ws.Cells.Item(29, 7).Value2 = 781242

$ws.Cells.Item(30, 3).Value2 = "45645018"
$ws.Cells.Item(30, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(30, 5).Value2 = "1805"
$ws.Cells.Item(30, 6).Value2 = 29509
$ws.Cells.Item(30, 7).Value2 = 781242

$ws.Cells.Item(31, 3).Value2 = "45645018"
$ws.Cells.Item(31, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(31, 5).Value2 = "1806"
$ws.Cells.Item(31, 6).Value2 = 29509
$ws.Cells.Item(31, 7).Value2 = 781242

$ws.Cells.Item(32, 3).Value2 = "45645018"
$ws.Cells.Item(32, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(32, 5).Value2 = "1807"
$ws.Cells.Item(32, 6).Value2 = 29509
$ws.Cells.Item(32, 7).Value2 = 781242

$ws.Cells.Item(33, 3).Value2 = "45645018"
$ws.Cells.Item(33, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(33, 5).Value2 = "1808"
$ws.Cells.Item(33, 6).Value2 = 29509
$ws.Cells.Item(33, 7).Value2 = 781242

$ws.Cells.Item(34, 3).Value2 = "45645018"
$ws.Cells.Item(34, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(34, 5).Value2 = "1809"
$ws.Cells.Item(34, 6).Value2 = 31249
$ws.Cells.Item(34, 7).Value2 = 781242

$ws.Cells.Item(35, 3).Value2 = "45645018"
$ws.Cells.Item(35, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(35, 5).Value2 = "1810"
$ws.Cells.Item(35, 6).Value2 = 31249
$ws.Cells.Item(35, 7).Value2 = 781242

$ws.Cells.Item(36, 3).Value2 = "45645018"
$ws.Cells.Item(36, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(36, 5).Value2 = "1811"
$ws.Cells.Item(36, 6).Value2 = 31249
$ws.Cells.Item(36, 7).Value2 = 781242

$ws.Cells.Item(37, 3).Value2 = "45645018"
$ws.Cells.Item(37, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(37, 5).Value2 = "1812"
$ws.Cells.Item(37, 6).Value2 = 31249
$ws.Cells.Item(37, 7).Value2 = 781242

$ws.Cells.Item(38, 3).Value2 = "45645018"
$ws.Cells.Item(38, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(38, 5).Value2 = "1901"
$ws.Cells.Item(38, 6).Value2 = 31249
$ws.Cells.Item(38, 7).Value2 = 781242

$ws.Cells.Item(39, 3).Value2 = "45645018"
$ws.Cells.Item(39, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(39, 5).Value2 = "1902"
$ws.Cells.Item(39, 6).Value2 = 31249
$ws.Cells.Item(39, 7).Value2 = 781242

$ws.Cells.Item(40, 3).Value2 = "45645018"
$ws.Cells.Item(40, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(40, 5).Value2 = "1903"
$ws.Cells.Item(40, 6).Value2 = 31249
$ws.Cells.Item(40, 7).Value2 = 781242

$ws.Cells.Item(41, 3).Value2 = "45645018"
$ws.Cells.Item(41, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(41, 5).Value2 = "1904"
$ws.Cells.Item(41, 6).Value2 = 31249
$ws.Cells.Item(41, 7).Value2 = 781242

$ws.Cells.Item(42, 3).Value2 = "45645018"
$ws.Cells.Item(42, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(42, 5).Value2 = "1905"
$ws.Cells.Item(42, 6).Value2 = 31249
$ws.Cells.Item(42, 7).Value2 = 781242

$ws.Cells.Item(43, 3).Value2 = "45645018"
$ws.Cells.Item(43, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(43, 5).Value2 = "1906"
$ws.Cells.Item(43, 6).Value2 = 31249
$ws.Cells.Item(43, 7).Value2 = 781242

$ws.Cells.Item(44, 3).Value2 = "45645018"
$ws.Cells.Item(44, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(44, 5).Value2 = "1907"
$ws.Cells.Item(44, 6).Value2 = 31249
$ws.Cells.Item(44, 7).Value2 = 781242

$ws.Cells.Item(45, 3).Value2 = "45645018"
$ws.Cells.Item(45, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(45, 5).Value2 = "1908"
$ws.Cells.Item(45, 6).Value2 = 31249
$ws.Cells.Item(45, 7).Value2 = 781242

$ws.Cells.Item(46, 3).Value2 = "45645018"
$ws.Cells.Item(46, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(46, 5).Value2 = "1909"
$ws.Cells.Item(46, 6).Value2 = 31249
$ws.Cells.Item(46, 7).Value2 = 781242

$ws.Cells.Item(47, 3).Value2 = "45645018"
$ws.Cells.Item(47, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(47, 5).Value2 = "1910"
$ws.Cells.Item(47, 6).Value2 = 31249
$ws.Cells.Item(47, 7).Value2 = 781242

$ws.Cells.Item(48, 3).Value2 = "45645018"
$ws.Cells.Item(48, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(48, 5).Value2 = "1911"
$ws.Cells.Item(48, 6).Value2 = 31249
$ws.Cells.Item(48, 7).Value2 = 781242

$ws.Cells.Item(49, 3).Value2 = "45645018"
$ws.Cells.Item(49, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(49, 5).Value2 = "1912"
$ws.Cells.Item(49, 6).Value2 = 31249
$ws.Cells.Item(49, 7).Value2 = 781242

$ws.Cells.Item(50, 3).Value2 = "45645018"
$ws.Cells.Item(50, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(50, 5).Value2 = "2001"
$ws.Cells.Item(50, 6).Value2 = 31249
$ws.Cells.Item(50, 7).Value2 = 781242

$ws.Cells.Item(51, 3).Value2 = "45645018"
$ws.Cells.Item(51, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(51, 5).Value2 = "2002"
$ws.Cells.Item(51, 6).Value2 = 31249
$ws.Cells.Item(51, 7).Value2 = 781242

$ws.Cells.Item(52, 3).Value2 = "45645018"
$ws.Cells.Item(52, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(52, 5).Value2 = "2003"
$ws.Cells.Item(52, 6).Value2 = 31249
$ws.Cells.Item(52, 7).Value2 = 781242

$ws.Cells.Item(53, 3).Value2 = "45645018"
$ws.Cells.Item(53, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(53, 5).Value2 = "2004"
$ws.Cells.Item(53, 6).Value2 = 31249
$ws.Cells.Item(53, 7).Value2 = 781242

$ws.Cells.Item(54, 3).Value2 = "45645018"
$ws.Cells.Item(54, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(54, 5).Value2 = "2005"
$ws.Cells.Item(54, 6).Value2 = 31249
$ws.Cells.Item(54, 7).Value2 = 781242

$ws.Cells.Item(55, 3).Value2 = "45645018"
$ws.Cells.Item(55, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(55, 5).Value2 = "2006"
$ws.Cells.Item(55, 6).Value2 = 31249
$ws.Cells.Item(55, 7).Value2 = 781242

$ws.Cells.Item(56, 3).Value2 = "45645018"
$ws.Cells.Item(56, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(56, 5).Value2 = "2007"
$ws.Cells.Item(56, 6).Value2 = 31249
$ws.Cells.Item(56, 7).Value2 = 781242

$ws.Cells.Item(57, 3).Value2 = "45645018"
$ws.Cells.Item(57, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(57, 5).Value2 = "2008"
$ws.Cells.Item(57, 6).Value2 = 31249
$ws.Cells.Item(57, 7).Value2 = 781242

$ws.Cells.Item(58, 3).Value2 = "45645018"
$ws.Cells.Item(58, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(58, 5).Value2 = "2009"
$ws.Cells.Item(58, 6).Value2 = 31249
$ws.Cells.Item(58, 7).Value2 = 781242

$ws.Cells.Item(59, 3).Value2 = "45645018"
$ws.Cells.Item(59, 4).Value2 = "MILENA JOHANA ANILLO CARO"
$ws.Cells.Item(59, 5).Value2 = "2010"
$ws.Cells.Item(59, 6).Value2 = 26041
$ws.Cells.Item(59, 7).Value2 = 781242
